# Replace the XLL-quote formulas in row 2 (bank name lookups) with their
# static cached results. This mirrors pasting the computed values back over
# the add-in formulas: the `_xll.xlquoteName(..., "FDIC")` formulas in
# C2:F2 are removed and the cells become plain (shared-string) literals
# holding the bank names that were previously just the cached <v> results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = "JPMORGAN CHASE BANK, NATIONAL ASSOCIATION"
$ws.Range("D2").Value2 = "BANK OF AMERICA, NATIONAL ASSOCIATION"
$ws.Range("E2").Value2 = "CITIBANK, NATIONAL ASSOCIATION"
$ws.Range("F2").Value2 = "WELLS FARGO BANK, NATIONAL ASSOCIATION"
